$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: M1_PH, CM2_PH, CMN3_PH, CMN4_PH
$ws.Range("I18").Value = -0.06439574632716669
$ws.Range("J18").Value = 0.2820221037857273
$ws.Range("K18").Value = -0.05024023413823025
$ws.Range("L18").Value = 2.288968724969885

# Row 19: M1_PH, CM2_PH, CMN3_PH, CMN4_PH
$ws.Range("I19").Value = 0.3127488043338543
$ws.Range("J19").Value = 0.5765091068230044
$ws.Range("K19").Value = 0.01146066641156172
$ws.Range("L19").Value = 1.873828640831355
